# The "Dashboard" SmartArt diagram (on the slide with the "Dashboard"
# title) has a node whose text is "ETL ". Rename it to
# "De-normalized data " — editing a SmartArt node's text through
# Shape.SmartArt.AllNodes(...).TextFrame2.TextRange.Text updates both the
# diagram's data part (dgm:t) and its cached rendered drawing part
# (dsp:txBody) in one go.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasSmartArt) {
            $sa = $sh.SmartArt
            for ($ni = 1; $ni -le $sa.AllNodes.Count; $ni++) {
                $node = $sa.AllNodes.Item($ni)
                $tr = $node.TextFrame2.TextRange
                if ($tr.Text -eq "ETL ") {
                    $tr.Text = "De-normalized data "
                }
            }
        }
    }
}
